$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in T2
$ws.Range("T2").Value = 314539

# Update the active selection to T3
$ws.Range("T3").Select()
